$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.893.29'
$ws.Range("E2").Value = '  -0.96%  '
$ws.Range("D3").Value = '2.356.01'
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("E5").Value = '  -2.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.07'
$ws.Range("E6").Value = '  -1.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.42'
$ws.Range("E7").Value = '  -4.03%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +2.38%  '
$ws.Range("E10").Value = '  -2.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '59.36'
$ws.Range("E11").Value = '  +2.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '33.43'
$ws.Range("E12").Value = '  +4.48%  '
$ws.Range("E13").Value = '  -1.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.109'
$ws.Range("E14").Value = '  +0.07%  '
$ws.Range("D15").Value = '2.705.92'
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.46'
$ws.Range("E16").Value = '  -3.84%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.909'
$ws.Range("E17").Value = '  -1.56%  '
$ws.Range("D18").Value = '2.357.31'
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").Value = '43.807.09'
$ws.Range("E19").Value = '  -1.14%  '
$ws.Range("E20").Value = '  -0.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.70'
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '77.60'
$ws.Range("E22").Value = '  -1.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '256.61'
$ws.Range("E23").Value = '  -0.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.95'
$ws.Range("E24").Value = '  +13.22%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.79'
$ws.Range("E26").Value = '  +1.87%  '
$ws.Range("E27").Value = '  -2.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.60'
$ws.Range("E28").Value = '  -2.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.88'
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '177.44'
$ws.Range("E30").Value = '  +1.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.19'
$ws.Range("E31").Value = '  -5.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.128'
$ws.Range("E32").Value = '  -1.02%  '
$ws.Range("E33").Value = '  +0.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0756'
$ws.Range("E34").Value = '  -0.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.20'
$ws.Range("E35").Value = '  -4.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.50'
$ws.Range("E36").Value = '  +2.25%  '
$ws.Range("E37").Value = '  -1.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.45'
$ws.Range("E38").Value = '  -2.48%  '
$ws.Range("E39").Value = '  -4.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0277'
$ws.Range("E40").Value = '  -0.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '67.83'
$ws.Range("E41").Value = '  +26.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.16'
$ws.Range("E42").Value = '  +15.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.110'
$ws.Range("E43").Value = '  +9.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.31'
$ws.Range("E44").Value = '  +2.03%  '
$ws.Range("E45").Value = '  +3.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.97'
$ws.Range("E46").Value = '  -1.02%  '
$ws.Range("E47").Value = '  -0.38%  '
$ws.Range("E48").Value = '  -0.12%  '
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '99.88'
$ws.Range("E50").Value = '  -2.09%  '
$ws.Range("E51").Value = '  -5.03%  '
